# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across all eight crafting-job sheets with latest market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 75279.85000000001
$ws.Range("I9").Value = 83597.44500000001
$ws.Range("K9").Value = 83597.44500000001
$ws.Range("M9").Value = -83428.44500000001
$ws.Range("H86").Value = 4445.2856
$ws.Range("I86").Value = 2276.2727
$ws.Range("K86").Value = 2276.2727
$ws.Range("M86").Value = -1153.2727
$ws.Range("H89").Value = 4445.2856
$ws.Range("I89").Value = 2276.2727
$ws.Range("K89").Value = 11381.3635
$ws.Range("M89").Value = -5765.363499999999
$ws.Range("H121").Value = 1040
$ws.Range("J121").Value = 1040
$ws.Range("L121").Value = 3120
$ws.Range("N121").Value = -6614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3213.0173
$ws.Range("I32").Value = 3095.3655
$ws.Range("J32").Value = 4232.6665
$ws.Range("K32").Value = 3095.3655
$ws.Range("L32").Value = 4232.6665
$ws.Range("M32").Value = -2808.3655
$ws.Range("N32").Value = -4806.6665
$ws.Range("H45").Value = 3483.25
$ws.Range("I45").Value = 1710.4286
$ws.Range("J45").Value = 5965.2
$ws.Range("K45").Value = 1710.4286
$ws.Range("L45").Value = 5965.2
$ws.Range("M45").Value = -1333.4286
$ws.Range("N45").Value = -6719.2
$ws.Range("H46").Value = 18557.4
$ws.Range("I46").Value = 10698.5
$ws.Range("J46").Value = 21415.182
$ws.Range("K46").Value = 10698.5
$ws.Range("L46").Value = 21415.182
$ws.Range("M46").Value = -10379.5
$ws.Range("N46").Value = -22053.182
$ws.Range("H63").Value = 1684.2
$ws.Range("I63").Value = 1621.5
$ws.Range("J63").Value = 1935
$ws.Range("K63").Value = 1621.5
$ws.Range("L63").Value = 1935
$ws.Range("M63").Value = -935.5
$ws.Range("N63").Value = -3307
$ws.Range("H66").Value = 1684.2
$ws.Range("I66").Value = 1621.5
$ws.Range("J66").Value = 1935
$ws.Range("K66").Value = 8107.5
$ws.Range("L66").Value = 9675
$ws.Range("M66").Value = -4675.5
$ws.Range("N66").Value = -16539
$ws.Range("H102").Value = 4488.7
$ws.Range("I102").Value = 3422.5715
$ws.Range("J102").Value = 6976.3335
$ws.Range("K102").Value = 3422.5715
$ws.Range("L102").Value = 6976.3335
$ws.Range("M102").Value = -1800.5715
$ws.Range("N102").Value = -10220.3335
$ws.Range("H122").Value = 4351.3
$ws.Range("I122").Value = 4448.143
$ws.Range("J122").Value = 4299.154
$ws.Range("K122").Value = 13344.429
$ws.Range("L122").Value = 12897.462
$ws.Range("M122").Value = -10894.429
$ws.Range("N122").Value = -17797.462
$ws.Range("H141").Value = 119764.5
$ws.Range("J141").Value = 119764.5
$ws.Range("L141").Value = 119764.5
$ws.Range("N141").Value = -130124.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4130.8
$ws.Range("I99").Value = 3977.5
$ws.Range("K99").Value = 3977.5
$ws.Range("M99").Value = -2479.5
$ws.Range("H102").Value = 36665.332
$ws.Range("I102").Value = 4998.5
$ws.Range("K102").Value = 4998.5
$ws.Range("M102").Value = -1753.5
$ws.Range("H105").Value = 1636186.2
$ws.Range("I105").Value = 2286771
$ws.Range("J105").Value = 9724.75
$ws.Range("K105").Value = 2286771
$ws.Range("L105").Value = 9724.75
$ws.Range("M105").Value = -2285024
$ws.Range("N105").Value = -13218.75
$ws.Range("H110").Value = 97673
$ws.Range("J110").Value = 97673
$ws.Range("L110").Value = 97673
$ws.Range("N110").Value = -105853
$ws.Range("H137").Value = 99057.336
$ws.Range("J137").Value = 99057.336
$ws.Range("L137").Value = 99057.336
$ws.Range("N137").Value = -109257.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 12988
$ws.Range("I86").Value = 13597.777
$ws.Range("K86").Value = 13597.777
$ws.Range("M86").Value = -12474.777
$ws.Range("H89").Value = 12988
$ws.Range("I89").Value = 13597.777
$ws.Range("K89").Value = 67988.88499999999
$ws.Range("M89").Value = -62372.88499999999
$ws.Range("H94").Value = 1954.5454
$ws.Range("I94").Value = 1683
$ws.Range("K94").Value = 1683
$ws.Range("M94").Value = -1232
$ws.Range("H99").Value = 37066
$ws.Range("I99").Value = 11106.571
$ws.Range("K99").Value = 11106.571
$ws.Range("M99").Value = -9608.571
$ws.Range("H107").Value = 4131.3335
$ws.Range("I107").Value = 1197
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 1197
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 723
$ws.Range("N107").Value = -13840
$ws.Range("H122").Value = 4552.5557
$ws.Range("J122").Value = 4161.6665
$ws.Range("L122").Value = 12484.9995
$ws.Range("N122").Value = -17384.9995
$ws.Range("H126").Value = 37066
$ws.Range("I126").Value = 11106.571
$ws.Range("K126").Value = 33319.713
$ws.Range("M126").Value = -30849.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 10338.866
$ws.Range("I138").Value = 9263.462
$ws.Range("J138").Value = 17329
$ws.Range("K138").Value = 27790.386
$ws.Range("L138").Value = 51987
$ws.Range("M138").Value = -22650.386
$ws.Range("N138").Value = -62267

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H97").Value = 761.7742
$ws.Range("I97").Value = 778.375
$ws.Range("J97").Value = 704.8570999999999
$ws.Range("K97").Value = 778.375
$ws.Range("L97").Value = 704.8570999999999
$ws.Range("M97").Value = -282.375
$ws.Range("N97").Value = -1696.8571
$ws.Range("H107").Value = 668.5
$ws.Range("I107").Value = 421.75
$ws.Range("K107").Value = 421.75
$ws.Range("M107").Value = 1498.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1699.6666
$ws.Range("I40").Value = 1949.5
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 1949.5
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -1813.5
$ws.Range("N40").Value = -1472
$ws.Range("H74").Value = 81666.336
$ws.Range("I74").Value = 81666.336
$ws.Range("K74").Value = 81666.336
$ws.Range("M74").Value = -80668.336
$ws.Range("H77").Value = 81666.336
$ws.Range("I77").Value = 81666.336
$ws.Range("K77").Value = 244999.008
$ws.Range("M77").Value = -240007.008
$ws.Range("H93").Value = 1986332.5
$ws.Range("I93").Value = 945.8261
$ws.Range("K93").Value = 945.8261
$ws.Range("M93").Value = 302.1739
$ws.Range("H122").Value = 3439.9246
$ws.Range("J122").Value = 4539.8
$ws.Range("L122").Value = 13619.4
$ws.Range("N122").Value = -18519.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 19987
$ws.Range("I18").Value = 19991.25
$ws.Range("J18").Value = 19970
$ws.Range("K18").Value = 19991.25
$ws.Range("L18").Value = 19970
$ws.Range("M18").Value = -19818.25
$ws.Range("N18").Value = -20316
$ws.Range("H34").Value = 62499.5
$ws.Range("J34").Value = 62499.5
$ws.Range("L34").Value = 62499.5
$ws.Range("N34").Value = -62905.5
$ws.Range("H54").Value = 35666.668
$ws.Range("J54").Value = 41000
$ws.Range("L54").Value = 41000
$ws.Range("N54").Value = -42040
$ws.Range("H81").Value = 2179
$ws.Range("I81").Value = 1895
$ws.Range("J81").Value = 2250
$ws.Range("K81").Value = 3790
$ws.Range("L81").Value = 4500
$ws.Range("M81").Value = -2729
$ws.Range("N81").Value = -6622
$ws.Range("H84").Value = 2179
$ws.Range("I84").Value = 1895
$ws.Range("J84").Value = 2250
$ws.Range("K84").Value = 18950
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = -13646
$ws.Range("N84").Value = -33108
$ws.Range("H107").Value = 2661.182
$ws.Range("I107").Value = 1468.9333
$ws.Range("K107").Value = 4406.7999
$ws.Range("M107").Value = -2486.7999
$ws.Range("H132").Value = 558530.0600000001
$ws.Range("I132").Value = 3421
$ws.Range("J132").Value = 1668748.1
$ws.Range("K132").Value = 10263
$ws.Range("L132").Value = 5006244.300000001
$ws.Range("M132").Value = -7733
$ws.Range("N132").Value = -5011304.300000001
$ws.Range("H136").Value = 529449.2
$ws.Range("I136").Value = 3158.1875
$ws.Range("J136").Value = 3336334.8
$ws.Range("K136").Value = 9474.5625
$ws.Range("L136").Value = 10009004.4
$ws.Range("M136").Value = -6924.5625
$ws.Range("N136").Value = -10014104.4
